$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data values (columns A:AH) for rows 2-5 -- replaces the previous
# 5 data rows with a fresh 1000-sample batch (new accuracy run), and the
# old trailing row 6 is dropped entirely (handled further below).
$row2 = @(45131.50694444445, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
$row3 = @(45131.51388888889, 5.808, 4.397, 0.02, 13.479, 11.83, 4.217, 12.812, 7.127, 3.254, 5.01, 6.51, 5.92, 1.874, 4.551, 6.724, 3.745, 0.065, 0.655, 66.36, 12.601, 4.175, 8.75, 4.398, 1.108, 6.285, 3.836, 3.934, 4.252, 6.603, 0.721, 11.49, 2.799, 5.313)
$row4 = @(45131.52083333334, 17.035, 12.798, 0.408, 37.58, 31.429, 13.157, 47.979, 20.678, 9.375999999999999, 14.03, 15.692, 16.114, 4.496, 13.312, 19.282, 11.03, 0.151, 0.791, 198.358, 37.292, 12.272, 25.477, 13.194, 2.072, 24.055, 10.974, 9.943, 11.537, 16.4, 0.386, 43.904, 7.311, 15.401)
$row5 = @(45131.52777777778, 22.39, 16.81, 0.62, 49.1, 40.83, 17.43, 68.09, 27.15, 12.28, 18.34, 20.12, 20.97, 5.77, 17.5, 25.25, 14.51, 0.19, 0.87, 261.43, 49.07, 16.14, 33.4, 17.38, 2.54, 33.25, 14.37, 12.82, 15, 21.1, 0.26, 62.17, 9.449999999999999, 20.23)

$rows = @($row2, $row3, $row4, $row5)
for ($r = 0; $r -lt $rows.Length; $r++) {
    $rowArr = $rows[$r]
    for ($c = 0; $c -lt $rowArr.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $rowArr[$c]
    }
}

# A handful of columns (K, Q, AB, AD) got a hair narrower (custom width
# 8 -> 7). ColumnWidth = 6.17 round-trips to a stored width of exactly 7.
$ws.Columns.Item(11).ColumnWidth = 6.17
$ws.Columns.Item(17).ColumnWidth = 6.17
$ws.Columns.Item(28).ColumnWidth = 6.17
$ws.Columns.Item(30).ColumnWidth = 6.17

# Drop the old last data row (r=6); the sheet now ends at row 5 and the
# dimension shrinks from A1:AH6 to A1:AH5 automatically.
$ws.Rows.Item(6).Delete()